$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startA = 10002
$startB = 110021

for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 1).Value = $startA + $i
    $ws.Cells.Item($r, 2).Value = $startB + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

$ws.PageSetup.Orientation = 1

$ws.Rows("31:1048576").Select()
